# Updated symbol list on Thu Dec 15 04:33:24 UTC 2022 with GitHub Actions
# Refresh the Price (column D) and a couple of Volume(1h) (column E) values
# to the latest pull from coinranking.com. Price values are stored as text
# (they were inline strings in the source sheet), so we force text entry
# via NumberFormat "@" before assigning, which preserves exact formatting
# (trailing zeros, etc.) instead of Excel re-parsing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text entry (these numeric-looking strings must stay text, with
    # their exact formatting / trailing zeros preserved) by switching to
    # the Text number format just long enough to assign the value, then
    # restore the cell's style to Normal so no stray formatting is left
    # behind on the cell (matches the source which carries no style here).
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "264.76"
Set-TextValue "D3"  "22.48"
Set-TextValue "D4"  "6.278"
Set-TextValue "D5"  "0.06150"
Set-TextValue "D6"  "3.593"
Set-TextValue "D7"  "6.663"
Set-TextValue "D8"  "1.346"
Set-TextValue "D9"  "0.8300"
Set-TextValue "D10" "0.01358"
Set-TextValue "D11" "0.1594"
Set-TextValue "D12" "0.08244"
Set-TextValue "D13" "0.03428"
Set-TextValue "D14" "0.03140"
Set-TextValue "D15" "0.09241"
Set-TextValue "D16" "3.900"
Set-TextValue "D17" "0.001711"
Set-TextValue "D19" "0.006253"
Set-TextValue "D20" "0.005271"
Set-TextValue "D24" "2.266"
Set-TextValue "D25" "0.3340"
Set-TextValue "D26" "0.1237"
Set-TextValue "D27" "0.0002679"
Set-TextValue "D40" "0.04622"
Set-TextValue "D41" "0.006970"
Set-TextValue "D42" "0.1136"

Set-TextValue "D43" "0.003248"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.01192"
Set-TextValue "D45" "0.00006145"
Set-TextValue "D46" "0.00000000751"

Set-TextValue "D47" "0.6999"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.1935"
Set-TextValue "D49" "0.00002102"
